$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "letter_group" column (D),
# pushing it from D to G. The new columns will hold max / min_non_zero / min.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "max"
$ws.Range("E1").Value = "min_non_zero"
$ws.Range("F1").Value = "min"

# Data rows: max (D), min_non_zero (E), min (F)
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0

$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0

$ws.Range("D9").Value = 1520
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0

$ws.Range("D10").Value = 1296
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 0

$ws.Range("D11").Value = 467
$ws.Range("E11").Value = 467
$ws.Range("F11").Value = 0

$ws.Range("D12").Value = 424
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 0

$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 0

$ws.Range("D14").Value = 101
$ws.Range("E14").Value = 101
$ws.Range("F14").Value = 0

$ws.Range("D15").Value = 21
$ws.Range("E15").Value = 21
$ws.Range("F15").Value = 0

# Tiny floating point recalculation drift in the "sem" column (C), matching
# the target snapshot exactly.
$ws.Range("C2").Value = 0.03664513893725978
$ws.Range("C3").Value = 0.0215601090599315
$ws.Range("C4").Value = 0.009009009009009011
$ws.Range("C5").Value = 0.02729866226453556
$ws.Range("C6").Value = 0.009009009009009011
$ws.Range("C7").Value = 0.009009009009009011
$ws.Range("C9").Value = 22.76734200418761
$ws.Range("C11").Value = 4.207207207207208
$ws.Range("C12").Value = 4.714935718036188
$ws.Range("C13").Value = 0.06306306306306306
$ws.Range("C14").Value = 0.9099099099099099
